$wb = $excel.ActiveWorkbook

# --- Sheet: safety_orders ---
$ws1 = $wb.Worksheets.Item("safety_orders")

$ws1.Range("D2").Value = 568.745654736
$ws1.Range("E2").Value = 581.450606268
$ws1.Range("F2").Value = 587.26511233068
$ws1.Range("G2").Value = 3.153508901828317

$ws1.Range("D3").Value = 538.72743138816
$ws1.Range("E3").Value = 560.08901882808
$ws1.Range("F3").Value = 565.6899090163608
$ws1.Range("G3").Value = 4.766299910685001

$ws1.Range("D4").Value = 491.899002953365
$ws1.Range("E4").Value = 525.9940108907225
$ws1.Range("F4").Value = 531.2539509996298
$ws1.Range("G4").Value = 7.407935126357701

$ws1.Range("D5").Value = 418.846654619414
$ws1.Range("E5").Value = 472.4203327550683
$ws1.Range("F5").Value = 477.144536082619
# G5 unchanged: 12.21807587735019

$ws1.Range("D6").Value = 304.88499120142
$ws1.Range("E6").Value = 388.6526619782442
$ws1.Range("F6").Value = 392.5391885980266
$ws1.Range("G6").Value = 22.33005007975585

Write-Host "sheet1 done"

# --- Sheet: open_buy_orders ---
$ws2 = $wb.Worksheets.Item("open_buy_orders")
$ws2.Range("A2").Value = "OPMBXR-NPKGL-ELMYOV"
$ws2.Range("B2").Value = 600.09
$ws2.Rows.Item(3).Delete()

Write-Host "sheet2 done"

# --- Sheet: open_sell_orders ---
$ws3 = $wb.Worksheets.Item("open_sell_orders")
$ws3.Range("A2").Value = "OM25BV-KJA3H-7HLR4M"
$ws3.Range("A3").Value = "O5A2K4-OKJPP-KUKHM3"

Write-Host "sheet3 done"

